$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$text_D12 = @'
Background
id="Par1">A COVID-19 outbreak started in Wuhan, China, last December and now has become a global pandemic.

 The clinical information in caring of critically ill patients with COVID-19 needs to be shared timely, especially under the situations that there is still a largely ongoing spread of COVID-19 in many countries.


Methods
id="Par2">A multicenter prospective observational study investigated all the COVID-19 patients received in 19 ICUs of 16 hospitals in Wuhan, China, over 24 h between 8 AM February 2h and 8 AM February 27, 2020. The demographic information, clinical characteristics, vital signs, complications, laboratory values, and clinical managements of the patients were studied.


Results
id="Par3">A total of 226 patients were included.

 Their median (interquartile range, IQR) age was 64 (57–70) years, and 139 (61.5%) patients were male.

 The duration from the date of ICU admission to the study date was 11 (5–17) days, and the duration from onset of symptoms to the study date was 31 (24–36) days.

 Among all the patients, 155 (68.6%) had at least one coexisting disease, and their sequential organ failure assessment score was 4 (2–8).

 Organ function damages were found in most of the patients: ARDS in 161 (71.2%) patients, septic shock in 34 (15.0%) patients, acute kidney injury occurred in 57 (25.2%) patients, cardiac injury in 61 (27.0%) patients, and lymphocytopenia in 160 (70.8%) patients.

 Of all the studied patients, 85 (37.6%) received invasive mechanical ventilation, including 14 (6.2%) treated with extracorporeal membrane oxygenation (ECMO) at the same time, 20 (8.8%) received noninvasive mechanical ventilation, and 24 (10.6%) received continuous renal replacement therapy.

 By April 9, 2020, 87 (38.5%) patients were deceased and 15 (6.7%) were still in the hospital.


Conclusions
id="Par4">Critically ill patients with COVID-19 are associated with a higher risk of severe complications and need to receive an intensive level of treatments.

 COVID-19 poses a great strain on critical care resources in hospitals.


Trial registration
id="Par5">Chinese Clinical Trial Registry, ChiCTR2000030164. Registered on February 24, 2020, http://www.

chictr.

org.

cn/edit.

aspx?pid=49983&amp;htm=4

'@
$ws.Range("D12").Value = $text_D12

$text_E12 = @'
[Yuan%Yu%NULL%4,  Dan%Xu%NULL%0,  Shouzhi%Fu%NULL%0,  Jun%Zhang%NULL%4,  Xiaobo%Yang%NULL%0,  Liang%Xu%NULL%4,  Jiqian%Xu%NULL%0,  Yongran%Wu%NULL%0,  Chaolin%Huang%NULL%4,  Yaqi%Ouyang%NULL%0,  Luyu%Yang%NULL%0,  Minghao%Fang%NULL%0,  Hongwen%Xiao%NULL%4,  Jing%Ma%NULL%0,  Wei%Zhu%NULL%4,  Song%Hu%NULL%4,  Quan%Hu%NULL%4,  Daoyin%Ding%NULL%4,  Ming%Hu%NULL%0,  Guochao%Zhu%NULL%0,  Weijiang%Xu%NULL%4,  Jun%Guo%NULL%4,  Jinglong%Xu%NULL%4,  Haitao%Yuan%NULL%4,  Bin%Zhang%NULL%5,  Zhui%Yu%yuzhui@whu.edu.cn%0,  Dechang%Chen%icudechangchen@163.com%4,  Shiying%Yuan%yuan_shiying@163.com%0,  You%Shang%you_shanghust@163.com%0]
'@
$ws.Range("E12").Value = $text_E12

$text_D13 = @'
Purpose
id="Par1">The coronavirus disease 2019 (COVID-19) outbreak has become a global public health concern; however, relatively few detailed reports of related cardiac injury are available.

 The aims of this study were to compare the clinical and echocardiographic characteristics of inpatients in the intensive-care unit (ICU) and non-ICU patients.


Methods
id="Par2">We recruited 416 patients diagnosed with COVID-19 and divided them into two groups: ICU (n = 35) and non-ICU (n = 381).

 Medical histories, laboratory findings, and echocardiography data were compared.


Results
id="Par3">The levels of myocardial injury markers in ICU vs non-ICU patients were as follows: troponin I (0.029 ng/mL [0.007–0.063] vs 0.006 ng/mL [0.006–0.006]) and myoglobin (65.45 μg/L [39.77–130.57] vs 37.00 μg/L [26.40–53.54]).

 Echocardiographic findings included ventricular wall thickening (12 [39%] vs 1 [4%]), pulmonary hypertension (9 [29%] vs 0 [0%]), and reduced left-ventricular ejection fraction (5 [16%] vs 0 [0%]).

 Overall, 10% of the ICU patients presented with right heart enlargement, thickened right-ventricular wall, decreased right heart function, and pericardial effusion.

 Cardiac complications were more common in ICU patients, including acute cardiac injury (21 [60%] vs 13 [3%]) (including 2 cases of fulminant myocarditis), atrial or ventricular tachyarrhythmia (3 [9%] vs 3 [1%]), and acute heart failure (5 [14%] vs 0 [0%]).


Conclusion
id="Par4">Myocardial injury marker elevation, ventricular wall thickening, pulmonary artery hypertension, and cardiac complications including acute myocardial injury, arrhythmia, and acute heart failure are more common in ICU patients with COVID-19. Cardiac injury in COVID-19 patients may be related more to the systemic response after infection rather than direct damage by coronavirus.



'@
$ws.Range("D13").Value = $text_D13

$text_E13 = @'
[Jia-Hui%Zeng%zengjiahui1993@163.com%1,  Wei-Bo%Wu%13923451469@139.com%2,  Wei-Bo%Wu%13923451469@139.com%0,  Jiu-Xin%Qu%qujiuxin@163.com%1,  Yao%Wang%68848030@qq.com%1,  Chang-Feng%Dong%dongchangfeng7@163.com%1,  Yong-Fang%Luo%luoyongfang2005@foxmail.com%1,  Dan%Zhou%506057507@qq.com%1,  Wen-Xia%Feng%1484177993@qq.com%1,  Cheng%Feng%chaosheng-01@szsy.sustech.edu.cn%0]
'@
$ws.Range("E13").Value = $text_E13

$text_D14 = @'
Background
id="Par1">Since December 2019, the 2019 coronavirus disease (COVID-19) has expanded to cause a worldwide outbreak that more than 600,000 people infected and tens of thousands died.

 To date, the clinical characteristics of COVID-19 patients in the non-Wuhan areas of Hubei Province in China have not been described.


Methods
id="Par2">We retrospectively analyzed the clinical characteristics and treatment progress of 91 patients diagnosed with COVID-19 in Jingzhou Central Hospital.


Results
id="Par3">Of the 91 patients diagnosed with COVID-19, 30 cases (33.0%) were severe and two patients (2.2%) died.

 The severe disease group tended to be older (50.5 vs.

 42.0 years; p = 0.049) and have more chronic disease (40% vs.

 14.8%; p = 0.009) relative to mild disease group.

 Only 73.6% of the patients were quantitative polymerase chain reaction (qPCR)-positive on their first tests, while typical chest computed tomography images were obtained for each patient.

 The most common complaints were cough (n = 75; 82.4%), fever (n = 59; 64.8%), fatigue (n = 35; 38.5%), and diarrhea (n = 14; 15.4%).

 Non-respiratory injury was identified by elevated levels of aspartate aminotransferase (n = 18; 19.8%), creatinine (n = 5; 5.5%), and creatine kinase (n = 14; 15.4%) in laboratory tests.

 Twenty-eight cases (30.8%) suffered non-respiratory injury, including 50% of the critically ill patients and 21.3% of the mild patients.


Conclusions
id="Par4">Overall, the mortality rate of patients in Jingzhou was lower than that of Wuhan.

 Importantly, we found liver, kidney, digestive tract, and heart injuries in COVID-19 cases besides respiratory problems.

 Combining chest computed tomography images with the qPCR analysis of throat swab samples can improve the accuracy of COVID-19 diagnosis.



'@
$ws.Range("D14").Value = $text_D14

$text_E14 = @'
[Xin-Ying%Zhao%NULL%1,  Xuan-Xuan%Xu%NULL%1,  Hai-Sen%Yin%NULL%1,  Qin-Ming%Hu%NULL%1,  Tao%Xiong%NULL%1,  Yuan-Yan%Tang%NULL%1,  Ai-Ying%Yang%NULL%1,  Bao-Ping%Yu%yubp62@163.com%1,  Zhi-Ping%Huang%huangzhipingdoctor@163.com%1]
'@
$ws.Range("E14").Value = $text_E14

$text_D31 = @'
Background
id="Par1">Myocardial involvement in the course of coronavirus disease 2019 (COVID-19) pneumonia has been reported, though not fully characterized yet.

 The aim of the present study is to undertake a joint evaluation of hs-Troponin and natriuretic peptides (NP) in patients hospitalized for COVID-19 pneumonia.


Methods
id="Par2">In this multicenter observational study, we analyzed data from n = 111 patients.

 Cardiac biomarkers subgroups were identified according to values beyond reference range.


Results
id="Par3">Increased hs-Troponin and NP were found in 38 and 56% of the cases, respectively.

 As compared to those with normal cardiac biomarkers, these patients were older, had higher prevalence of cardiovascular diseases (CVD) and had more severe COVID-19 pneumonia by higher CRP and d-dimer and lower PaO2/FIO2. Two-dimensional echocardiography performed in a subset of patients (n = 24) showed significantly reduced left ventricular ejection fraction in patients with elevated NP (p = 0.02), whereas right ventricular systolic function (tricuspid annular plane systolic excursion) was significantly reduced both in patients with high hs-Troponin and NP (p = 0.022 and p = 0.03, respectively).

 Both hs-Troponin and NP were higher in patients with in-hospital mortality (p = 0.001 and p = 0.002, respectively).

 On multivariable analysis, independent associations were found of hs-Troponin with age, PaO2/FIO2 and d-dimer (B = 0.419, p = 0.001; B = − 0.212, p = 0.013; and B = 0.179, p = 0.037, respectively) and of NP with age and previous CVD (B = 0.480, p &lt; 0.001; and B = 0.253, p = 0.001, respectively).


Conclusions
id="Par4">Myocardial involvement at admission is common in COVID-19 pneumonia.

 Independent associations of hs-Troponin with markers of disease severity and of NP with underlying CVD might point toward existing different mechanisms leading to their elevation in this setting.


Electronic supplementary material
The online version of this article (10.1007/s11739-020-02498-7) contains supplementary material, which is available to authorized users.



'@
$ws.Range("D31").Value = $text_D31

$text_E31 = @'
[Luca%Arcari%luca_arcari@outlook.it%1,  Michelangelo%Luciani%NULL%2,  Michelangelo%Luciani%NULL%0,  Luca%Cacciotti%NULL%1,  Maria Beatrice%Musumeci%NULL%1,  Valerio%Spuntarelli%NULL%1,  Eleonora%Pistella%NULL%1,  Dario%Martolini%NULL%1,  Daniele%Manzo%NULL%1,  Mariateresa%Pucci%NULL%1,  Claudio%Marone%NULL%1,  Serena%Melandri%NULL%1,  Gerardo%Ansalone%NULL%1,  Claudio%Santini%NULL%1,  Paolo%Martelletti%NULL%1,  Massimo%Volpe%NULL%1,  Luciano%De Biase%NULL%1]
'@
$ws.Range("E31").Value = $text_E31

$text_D32 = @'
Purpose
id="Par1">Coronavirus disease 2019 (COVID-19) is creating an unprecedented healthcare crisis.

 Understanding the determinants of mortality is crucial to optimise intensive care unit (ICU) resource use and to identify targets for improving survival.


Methods
id="Par2">In a multicentre retrospective study, we included 379 COVID-19 patients admitted to four ICUs between 20 February and 24 April 2020 and categorised according to time from disease onset to ICU admission.

 A Cox proportional-hazards model identified factors associated with 28-day mortality.


Results
id="Par3">Median age was 66 years (53–68) and 292 (77%) were men.

 The main comorbidities included obesity and overweight (67%), hypertension (49.6%) and diabetes (30.1%).

 Median time from disease onset (i.

e.

, viral symptoms) to ICU admission was 8 (6–11) days (missing for three); 161 (42.5%) patients were admitted within a week of disease onset, 173 (45.6%) between 8 and 14 days, and 42 (11.1%) &gt; 14 days after disease onset; day 28 mortality was 26.4% (22–31) and decreased as time from disease onset to ICU admission increased, from 37 to 21% and 12%, respectively.

 Patients admitted within the first week had higher SOFA scores, more often had thrombocytopenia or acute kidney injury, had more limited radiographic involvement, and had significantly higher blood IL-6 levels.

 Age, COPD, immunocompromised status, time from disease onset, troponin concentration, and acute kidney injury were independently associated with mortality.


Conclusion
id="Par4">The excess mortality in patients admitted within a week of disease onset reflected greater non-respiratory severity.

 Therapeutic interventions against SARS-CoV-2 might impact different clinical endpoints according to time since disease onset.


Electronic supplementary material
The online version of this article (10.1007/s00134-020-06202-3) contains supplementary material, which is available to authorized users.



'@
$ws.Range("D32").Value = $text_D32

$text_E32 = @'
[Elie%Azoulay%elie.azoulay@aphp.fr%1,  Muriel%Fartoukh%NULL%2,  Muriel%Fartoukh%NULL%0,  Michael%Darmon%NULL%1,  Guillaume%Géri%NULL%1,  Guillaume%Voiriot%NULL%1,  Thibault%Dupont%NULL%1,  Lara%Zafrani%NULL%1,  Lola%Girodias%NULL%1,  Vincent%Labbé%NULL%1,  Martin%Dres%NULL%1,  Alexandra%Beurton%NULL%1,  Antoine%Vieillard-Baron%NULL%1,  Alexandre%Demoule%NULL%1]
'@
$ws.Range("E32").Value = $text_E32

$text_D34 = @'
Backgrounds
id="Par1">Patients at greatest risk of severe clinical conditions from coronavirus disease 2019 (COVID-19) and death are elderly and comorbid patients.

 Increased levels of cardiac troponins identify patients with poor outcome.

 The present study aimed to describe the clinical characteristics and outcomes of a cohort of Italian inpatients, admitted to a medical COVID-19 Unit, and to investigate the relative role of cardiac injury on in-hospital mortality.


Methods and results
id="Par2">We analyzed all consecutive patients with laboratory-confirmed COVID-19 referred to our dedicated medical Unit between February 26th and March 31st 2020. Patients’ clinical data including comorbidities, laboratory values, and outcomes were collected.

 Predictors of in-hospital mortality were investigated.

 A mediation analysis was performed to identify the potential mediators in the relationship between cardiac injury and mortality.

 A total of 109 COVID-19 inpatients (female 36%, median age 71 years) were included.

 During in-hospital stay, 20 patients (18%) died and, compared with survivors, these patients were older, had more comorbidities defined by Charlson comorbidity index ≥ 3(65% vs 24%, p = 0.001), and higher levels of high-sensitivity cardiac troponin I (Hs-cTnI), both at first evaluation and peak levels.

 A dose–response curve between Hs-cTnI and in-hospital mortality risk up to 200 ng/L was detected.

 Hs-cTnI, chronic kidney disease, and chronic coronary artery disease mediated most of the risk of in-hospital death, with Hs-cTnI mediating 25% of such effect.

 Smaller effects were observed for age, lactic dehydrogenase, and d-dimer.


Conclusions
id="Par3">In this cohort of elderly and comorbid COVID-19 patients, elevated Hs-cTnI levels were the most important and independent mediators of in-hospital mortality.


Electronic supplementary material
The online version of this article (10.1007/s11739-020-02495-w) contains supplementary material, which is available to authorized users.



'@
$ws.Range("D34").Value = $text_D34

$text_E34 = @'
[Alberto%Cipriani%alberto.cipriani@unipd.it%1,  Federico%Capone%NULL%2,  Federico%Capone%NULL%0,  Filippo%Donato%NULL%1,  Leonardo%Molinari%NULL%1,  Davide%Ceccato%NULL%1,  Alois%Saller%NULL%1,  Lorenzo%Previato%NULL%1,  Raffaele%Pesavento%NULL%1,  Cristiano%Sarais%NULL%1,  Paola%Fioretto%NULL%2,  Sabino%Iliceto%NULL%1,  Dario%Gregori%NULL%1,  Angelo%Avogaro%NULL%1,  Roberto%Vettor%NULL%3]
'@
$ws.Range("E34").Value = $text_E34
